$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1117604.5
$ws.Range("J17").Value = 1771539.8
$ws.Range("L17").Value = 5314619.4
$ws.Range("N17").Value = -5314955.4

$ws.Range("H21").Value = 15912.556
$ws.Range("I21").Value = 9905.666999999999
$ws.Range("J21").Value = 18916
$ws.Range("K21").Value = 9905.666999999999
$ws.Range("L21").Value = 18916
$ws.Range("M21").Value = -9437.666999999999
$ws.Range("N21").Value = -19852

$ws.Range("H23").Value = 15912.556
$ws.Range("I23").Value = 9905.666999999999
$ws.Range("J23").Value = 18916
$ws.Range("K23").Value = 9905.666999999999
$ws.Range("L23").Value = 18916
$ws.Range("M23").Value = -9671.666999999999
$ws.Range("N23").Value = -19384

$ws.Range("H29").Value = 303.2857
$ws.Range("I29").Value = 103.833336
$ws.Range("J29").Value = 1500
$ws.Range("K29").Value = 311.500008
$ws.Range("L29").Value = 4500
$ws.Range("M29").Value = -30.50000799999998
$ws.Range("N29").Value = -5062

$ws.Range("H51").Value = 1937.5
$ws.Range("I51").Value = 1833.3334
$ws.Range("K51").Value = 1833.3334
$ws.Range("M51").Value = -1349.3334

$ws.Range("H116").Value = 14951.25
$ws.Range("I116").Value = 22162
$ws.Range("J116").Value = 2933.3333
$ws.Range("K116").Value = 22162
$ws.Range("L116").Value = 2933.3333
$ws.Range("M116").Value = -18720
$ws.Range("N116").Value = -9817.3333

$ws.Range("H121").Value = 2143.75
$ws.Range("J121").Value = 3900
$ws.Range("L121").Value = 11700
$ws.Range("N121").Value = -15194

$ws.Range("H131").Value = 3178.7144
$ws.Range("I131").Value = 1899
$ws.Range("J131").Value = 4138.5
$ws.Range("K131").Value = 5697
$ws.Range("L131").Value = 12415.5
$ws.Range("M131").Value = -657
$ws.Range("N131").Value = -22495.5

$ws.Range("H141").Value = 2085.1453
$ws.Range("I141").Value = 727.6799999999999
$ws.Range("J141").Value = 7741.25
$ws.Range("K141").Value = 2183.04
$ws.Range("L141").Value = 23223.75
$ws.Range("M141").Value = 2996.96
$ws.Range("N141").Value = -33583.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 41000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 41000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 41000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -41228

$ws.Range("H32").Value = 517402.56
$ws.Range("I32").Value = 568853.1
$ws.Range("J32").Value = 22191
$ws.Range("K32").Value = 568853.1
$ws.Range("L32").Value = 22191
$ws.Range("M32").Value = -568566.1
$ws.Range("N32").Value = -22765

$ws.Range("H45").Value = 2908.5386
$ws.Range("I45").Value = 2946.4546
$ws.Range("K45").Value = 2946.4546
$ws.Range("M45").Value = -2569.4546

$ws.Range("H74").Value = 1371.9667
$ws.Range("I74").Value = 930.7368
$ws.Range("J74").Value = 2134.0908
$ws.Range("K74").Value = 930.7368
$ws.Range("L74").Value = 2134.0908
$ws.Range("M74").Value = -56.73680000000002
$ws.Range("N74").Value = -3882.0908

$ws.Range("H77").Value = 1371.9667
$ws.Range("I77").Value = 930.7368
$ws.Range("J77").Value = 2134.0908
$ws.Range("K77").Value = 4653.684
$ws.Range("L77").Value = 10670.454
$ws.Range("M77").Value = -285.6840000000002
$ws.Range("N77").Value = -19406.454

$ws.Range("H122").Value = 68948
$ws.Range("I122").Value = 127052.75
$ws.Range("J122").Value = 2542.5715
$ws.Range("K122").Value = 381158.25
$ws.Range("L122").Value = 7627.7145
$ws.Range("M122").Value = -378708.25
$ws.Range("N122").Value = -12527.7145

$ws.Range("H123").Value = 24428.428
$ws.Range("J123").Value = 24428.428
$ws.Range("L123").Value = 24428.428
$ws.Range("N123").Value = -34228.428

$ws.Range("H132").Value = 2147.182
$ws.Range("I132").Value = 1254.5818
$ws.Range("K132").Value = 3763.7454
$ws.Range("M132").Value = -1233.7454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1994.5834
$ws.Range("I99").Value = 1976.4286
$ws.Range("J99").Value = 2020
$ws.Range("K99").Value = 1976.4286
$ws.Range("L99").Value = 2020
$ws.Range("M99").Value = -478.4286
$ws.Range("N99").Value = -5016

$ws.Range("H124").Value = 50780
$ws.Range("J124").Value = 50780
$ws.Range("L124").Value = 50780
$ws.Range("N124").Value = -60600

$ws.Range("H134").Value = 2805.439
$ws.Range("I134").Value = 2408.923
$ws.Range("J134").Value = 3492.7334
$ws.Range("K134").Value = 7226.768999999999
$ws.Range("L134").Value = 10478.2002
$ws.Range("M134").Value = -4691.768999999999
$ws.Range("N134").Value = -15548.2002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H22").Value = 83333900
$ws.Range("I22").Value = 380
$ws.Range("K22").Value = 380
$ws.Range("M22").Value = -30

$ws.Range("H31").Value = 4991.52
$ws.Range("I31").Value = 1151.375
$ws.Range("J31").Value = 8536.27
$ws.Range("K31").Value = 1151.375
$ws.Range("L31").Value = 8536.27
$ws.Range("M31").Value = -856.375
$ws.Range("N31").Value = -9126.27

$ws.Range("H34").Value = 4991.52
$ws.Range("I34").Value = 1151.375
$ws.Range("J34").Value = 8536.27
$ws.Range("K34").Value = 1151.375
$ws.Range("L34").Value = 8536.27
$ws.Range("M34").Value = -949.375
$ws.Range("N34").Value = -8940.27

$ws.Range("H122").Value = 1980.48
$ws.Range("I122").Value = 1753.6666
$ws.Range("K122").Value = 5260.9998
$ws.Range("M122").Value = -2810.9998

$ws.Range("H132").Value = 44873604
$ws.Range("I132").Value = 55557220
$ws.Range("J132").Value = 20835468
$ws.Range("K132").Value = 166671660
$ws.Range("L132").Value = 62506404
$ws.Range("M132").Value = -166669130
$ws.Range("N132").Value = -62511464

$ws.Range("H141").Value = 189285.72
$ws.Range("J141").Value = 187500
$ws.Range("L141").Value = 187500
$ws.Range("N141").Value = -197860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 30000858
$ws.Range("I4").Value = 28000374
$ws.Range("J4").Value = 33335000
$ws.Range("K4").Value = 84001122
$ws.Range("L4").Value = 100005000
$ws.Range("M4").Value = -84001010
$ws.Range("N4").Value = -100005224

$ws.Range("H38").Value = 8342144.5
$ws.Range("I38").Value = 12513192
$ws.Range("J38").Value = 50
$ws.Range("K38").Value = 37539576
$ws.Range("L38").Value = 150
$ws.Range("M38").Value = -37539229
$ws.Range("N38").Value = -844

$ws.Range("H109").Value = 3836.5
$ws.Range("I109").Value = 890.875
$ws.Range("J109").Value = 6193
$ws.Range("K109").Value = 2672.625
$ws.Range("L109").Value = 18579
$ws.Range("M109").Value = -1632.625
$ws.Range("N109").Value = -20659

$ws.Range("H121").Value = 1393.4286
$ws.Range("I121").Value = 394.42856
$ws.Range("J121").Value = 1892.9286
$ws.Range("K121").Value = 1183.28568
$ws.Range("L121").Value = 5678.7858
$ws.Range("M121").Value = 126.71432
$ws.Range("N121").Value = -8298.7858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1166.6666
$ws.Range("I122").Value = 1250
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 3750
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -1300
$ws.Range("N122").Value = -7900

$ws.Range("H132").Value = 2985.566
$ws.Range("I132").Value = 2619.7805
$ws.Range("J132").Value = 4235.3335
$ws.Range("K132").Value = 7859.3415
$ws.Range("L132").Value = 12706.0005
$ws.Range("M132").Value = -5329.3415
$ws.Range("N132").Value = -17766.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9692.462
$ws.Range("I22").Value = 1999.875
$ws.Range("J22").Value = 22000.6
$ws.Range("K22").Value = 1999.875
$ws.Range("L22").Value = 22000.6
$ws.Range("M22").Value = -1704.875
$ws.Range("N22").Value = -22590.6

$ws.Range("H27").Value = 9692.462
$ws.Range("I27").Value = 1999.875
$ws.Range("J27").Value = 22000.6
$ws.Range("K27").Value = 1999.875
$ws.Range("L27").Value = 22000.6
$ws.Range("M27").Value = -1892.875
$ws.Range("N27").Value = -22214.6

$ws.Range("H46").Value = 2299.9312
$ws.Range("I46").Value = 1747.92
$ws.Range("J46").Value = 5750
$ws.Range("K46").Value = 1747.92
$ws.Range("L46").Value = 5750
$ws.Range("M46").Value = -1559.92
$ws.Range("N46").Value = -6126

$ws.Range("H122").Value = 3340.12
$ws.Range("I122").Value = 2960.3
$ws.Range("J122").Value = 3593.3333
$ws.Range("K122").Value = 8880.900000000001
$ws.Range("L122").Value = 10779.9999
$ws.Range("M122").Value = -6430.900000000001
$ws.Range("N122").Value = -15679.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 250052500
$ws.Range("I4").Value = 1000000000
$ws.Range("J4").Value = 70003
$ws.Range("K4").Value = 1000000000
$ws.Range("L4").Value = 70003
$ws.Range("M4").Value = -999999887
$ws.Range("N4").Value = -70229

$ws.Range("H122").Value = 2123.4285
$ws.Range("I122").Value = 2100
$ws.Range("J122").Value = 2191.111
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 6573.333
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -11473.333

$ws.Range("H123").Value = 30785.8
$ws.Range("J123").Value = 46964.5
$ws.Range("L123").Value = 46964.5
$ws.Range("N123").Value = -56764.5

$ws.Range("H126").Value = 2342.8572
$ws.Range("I126").Value = 1200
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 3600
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -1130
$ws.Range("N126").Value = -13340

$ws.Range("H136").Value = 2240.55
$ws.Range("I136").Value = 1857.8478
$ws.Range("J136").Value = 3498
$ws.Range("K136").Value = 5573.5434
$ws.Range("L136").Value = 10494
$ws.Range("M136").Value = -3023.5434
$ws.Range("N136").Value = -15594

$ws.Range("H139").Value = 99715
$ws.Range("J139").Value = 99715
$ws.Range("L139").Value = 99715
$ws.Range("N139").Value = -109995
